# Add a new "StatQuery" column (B) for stat-bar validation, shifting the
# existing dbExcel/WebExcel file-name columns one slot to the right.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B, shifting old B (dbExcel) and C (WebExcel) right.
$ws.Range("B1").EntireColumn.Insert()

# New column B: stat query text first (so it lands at sharedStrings index 6),
# then the "StatQuery" header (landing at index 7) - matches author's save order.
$ws.Range("B2").Value = "MATCH (t:clinical_trial)<--(a:arm)<--(c:case)<--(s:specimen)<--(:assignment_report) WITH DISTINCT c AS c, t ,a, s WHERE t.clinical_trial_designation IN ['NCI-MATCH'] OPTIONAL MATCH (s)<-[*]-(f:file) RETURN count(DISTINCT(f)) as number_of_files , count(DISTINCT(c.case_id)) as number_of_cases , count(DISTINCT(t.clinical_trial_designation)) as number_of_trial"
$ws.Range("B1").Value = "StatQuery"

# Apply wrap-text style (same as A2) to the new B2 cell
$ws.Range("B2").WrapText = $true

# Match column B's width to column A (75.81640625 "characters"); columns A, C, D
# already retain their original widths across the insert and need no change.
# (75.0 is the closest input that round-trips to the nearest achievable width.)
$ws.Range("B:B").ColumnWidth = 75.0

# Update selection / view
$ws.Range("B2").Select()
